$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached "last printed/edited" date fields
#    (Handout Master & Notes Master date placeholders): 01/09/2017 -> 05/09/2017
# ---------------------------------------------------------------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "05/09/2017"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "05/09/2017"

# ---------------------------------------------------------------------
# 2) Slide 1 - presenter/subtitle text box: re-enter the "Koos Drost"
#    and "05-09-2017" lines so they are stored as plain, single-run
#    paragraphs.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subTr = $subtitle.TextFrame.TextRange
$para2 = $subTr.Paragraphs(2, 1)
$para2.Text = "Koos Drost"
$para3 = $subTr.Paragraphs(3, 1)
$para3.Text = "05-09-2017"

# ---------------------------------------------------------------------
# 3) Slide 2 - merge the closing-quote run into the "annotations" run
#    so the quote reads: annotations"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$quoteShape = $s2.Shapes.Item(2)
$quoteTr = $quoteShape.TextFrame.TextRange
$quoteFull = $quoteTr.Text
$annIdx = $quoteFull.IndexOf("annotations")
$annAndQuote = $quoteTr.Characters($annIdx + 1, $quoteFull.Length - $annIdx)
$annAndQuote.Text = "annotations" + [char]0x201D

# ---------------------------------------------------------------------
# 4) Slide 4 - "Lab 1" textbox: split the first line into "Lab 1 " plus
#    a hyperlinked GitHub URL, pointing at the workshop repository.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$labShape = $s4.Shapes.Item(2)
$labTr = $labShape.TextFrame.TextRange

# Expand the trailing "1 " run so it also carries the new URL text.
$growRange = $labTr.Characters(5, 2)
$growRange.Text = "1 https://github.com/koosdrost/workshop "

# Re-assert each run boundary to match the final run layout.
$runLab = $labTr.Characters(1, 4)
$runLab.Text = "Lab "

$runOne = $labTr.Characters(5, 2)
$runOne.Text = "1 "

$runHttps = $labTr.Characters(7, 5)
$runHttps.Text = "https"

$runSlashes = $labTr.Characters(12, 3)
$runSlashes.Text = "://"

$runGithub = $labTr.Characters(15, 29)
$runGithub.Text = "github.com/koosdrost/workshop"

$runTrailingSpace = $labTr.Characters(44, 1)
$runTrailingSpace.Text = " "

# Apply the hyperlink across the full URL ("https://github.com/koosdrost/workshop").
$urlRange = $labTr.Characters(7, 37)
$urlRange.ActionSettings(1).Hyperlink.Address = "https://github.com/koosdrost/workshop"
